# EMU <-> point conversion (PowerPoint COM Left/Top/Width/Height are in points)
$emuPerPt = 12700.0

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1) Shrink the "enter a question pls" input-box rectangle ("Rechthoek 5")
#    to make room for the new "send" button next to it.
#    off stays x=2950368 y=5743575 ; ext cx 4786311 -> 3512215 (cy 535781)
# ---------------------------------------------------------------------------
$inputBox = $s.Shapes.Item(3)
$inputBox.Width = 3512215 / $emuPerPt

# ---------------------------------------------------------------------------
# 2) Add the new "send" button ("Rechthoek 1").
#    Duplicate the input box so the new shape inherits the same bg1 fill and
#    p:style block (lnRef/fillRef/effectRef/fontRef), then reposition/resize
#    it, rename it and set its caption text.
# ---------------------------------------------------------------------------
$sendBtnRange = $inputBox.Duplicate()
$sendBtn = $sendBtnRange.Item(1)
$sendBtn.Name = "Rechthoek 1"
$sendBtn.Left = 6709719 / $emuPerPt
$sendBtn.Top = 5743575 / $emuPerPt
$sendBtn.Width = 1026959 / $emuPerPt
$sendBtn.Height = 535781 / $emuPerPt
$sendBtn.TextFrame.TextRange.Text = "send"
